$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 0) Create the size-8 "Calibri" font that Excel silently added to
#    styles.xml (referenced only by the worksheet's phoneticPr, which
#    this host does not expose). We stamp it on a throwaway cell far
#    outside any used range, then clear that cell so it leaves no
#    trace in the sheet data / dimension, only in the style table.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("Total")
$tmp = $wsTotal.Cells.Item(500, 500)
$tmp.Font.Size = 8
$tmp.Clear()

# ------------------------------------------------------------------
# 1) "Total Load" sheet view: scroll/selection change
# ------------------------------------------------------------------
$wsTotalLoad = $wb.Worksheets.Item("Total Load")
$wsTotalLoad.Activate()
$wsTotalLoad.Range("D12:E18").Select() | Out-Null

# ------------------------------------------------------------------
# 2) "Total-No-Diffusion" sheet: new rows 12-18 + sheet view
# ------------------------------------------------------------------
$wsTND = $wb.Worksheets.Item("Total-No-Diffusion")

$wsTND.Range("D12").Value = "C_0"
$wsTND.Range("D13").Value = "C_1"
$wsTND.Range("D14").Value = "C_2"
$wsTND.Range("D15").Value = "C_3"
$wsTND.Range("D16").Value = "C_4"
$wsTND.Range("D17").Value = "C_5"
$wsTND.Range("D18").Value = "C_6"

$wsTND.Range("E12:E18").FormulaArray = "=MMULT(MINVERSE(C3:I9),M3:M9)"

$wsTND.Rows.Item(12).RowHeight = 45.5
$wsTND.Rows.Item(13).RowHeight = 45.5
$wsTND.Rows.Item(14).RowHeight = 45.5
$wsTND.Rows.Item(15).RowHeight = 45.5
$wsTND.Rows.Item(16).RowHeight = 45.5
$wsTND.Rows.Item(17).RowHeight = 45.5
$wsTND.Rows.Item(18).RowHeight = 45.5

$wsTND.Activate()
$wsTND.Range("M3").Select() | Out-Null

# ------------------------------------------------------------------
# 3) "Total" sheet: fix formulas to include Placeholders!, add R4:U4,
#    add new rows 12-18, page setup, sheet view
# ------------------------------------------------------------------
$wsTotal.Range("C3").Formula = "='Total Load'!C3+Advection!C3+Diffusion!C3+Placeholders!C3"
$wsTotal.Range("D3").Formula = "='Total Load'!D3+Advection!D3+Diffusion!D3+Placeholders!D3"
$wsTotal.Range("E3").Formula = "='Total Load'!E3+Advection!E3+Diffusion!E3+Placeholders!E3"
$wsTotal.Range("F3").Formula = "='Total Load'!F3+Advection!F3+Diffusion!F3+Placeholders!F3"
$wsTotal.Range("G3").Formula = "='Total Load'!G3+Advection!G3+Diffusion!G3+Placeholders!G3"
$wsTotal.Range("H3").Formula = "='Total Load'!H3+Advection!H3+Diffusion!H3+Placeholders!H3"
$wsTotal.Range("I3").Formula = "='Total Load'!I3+Advection!I3+Diffusion!I3+Placeholders!I3"

$wsTotal.Range("C4").Formula = "='Total Load'!C4+Advection!C4+Diffusion!C4+Placeholders!C4"
$wsTotal.Range("D4").Formula = "='Total Load'!D4+Advection!D4+Diffusion!D4+Placeholders!D4"
$wsTotal.Range("E4").Formula = "='Total Load'!E4+Advection!E4+Diffusion!E4+Placeholders!E4"
$wsTotal.Range("F4").Formula = "='Total Load'!F4+Advection!F4+Diffusion!F4+Placeholders!F4"
$wsTotal.Range("G4").Formula = "='Total Load'!G4+Advection!G4+Diffusion!G4+Placeholders!G4"
$wsTotal.Range("H4").Formula = "='Total Load'!H4+Advection!H4+Diffusion!H4+Placeholders!H4"
$wsTotal.Range("I4").Formula = "='Total Load'!I4+Advection!I4+Diffusion!I4+Placeholders!I4"

$wsTotal.Range("C5").Formula = "='Total Load'!C5+Advection!C5+Diffusion!C5+Placeholders!C5"
$wsTotal.Range("D5").Formula = "='Total Load'!D5+Advection!D5+Diffusion!D5+Placeholders!D5"
$wsTotal.Range("E5").Formula = "='Total Load'!E5+Advection!E5+Diffusion!E5+Placeholders!E5"
$wsTotal.Range("F5").Formula = "='Total Load'!F5+Advection!F5+Diffusion!F5+Placeholders!F5"
$wsTotal.Range("G5").Formula = "='Total Load'!G5+Advection!G5+Diffusion!G5+Placeholders!G5"
$wsTotal.Range("H5").Formula = "='Total Load'!H5+Advection!H5+Diffusion!H5+Placeholders!H5"
$wsTotal.Range("I5").Formula = "='Total Load'!I5+Advection!I5+Diffusion!I5+Placeholders!I5"

$wsTotal.Range("C6").Formula = "='Total Load'!C6+Advection!C6+Diffusion!C6+Placeholders!C6"
$wsTotal.Range("D6").Formula = "='Total Load'!D6+Advection!D6+Diffusion!D6+Placeholders!D6"
$wsTotal.Range("E6").Formula = "='Total Load'!E6+Advection!E6+Diffusion!E6+Placeholders!E6"
$wsTotal.Range("F6").Formula = "='Total Load'!F6+Advection!F6+Diffusion!F6+Placeholders!F6"
$wsTotal.Range("G6").Formula = "='Total Load'!G6+Advection!G6+Diffusion!G6+Placeholders!G6"
$wsTotal.Range("H6").Formula = "='Total Load'!H6+Advection!H6+Diffusion!H6+Placeholders!H6"
$wsTotal.Range("I6").Formula = "='Total Load'!I6+Advection!I6+Diffusion!I6+Placeholders!I6"

$wsTotal.Range("C7").Formula = "='Total Load'!C7+Advection!C7+Diffusion!C7+Placeholders!C7"
$wsTotal.Range("D7").Formula = "='Total Load'!D7+Advection!D7+Diffusion!D7+Placeholders!D7"
$wsTotal.Range("E7").Formula = "='Total Load'!E7+Advection!E7+Diffusion!E7+Placeholders!E7"
$wsTotal.Range("F7").Formula = "='Total Load'!F7+Advection!F7+Diffusion!F7+Placeholders!F7"
$wsTotal.Range("G7").Formula = "='Total Load'!G7+Advection!G7+Diffusion!G7+Placeholders!G7"
$wsTotal.Range("H7").Formula = "='Total Load'!H7+Advection!H7+Diffusion!H7+Placeholders!H7"
$wsTotal.Range("I7").Formula = "='Total Load'!I7+Advection!I7+Diffusion!I7+Placeholders!I7"

$wsTotal.Range("C8").Formula = "='Total Load'!C8+Advection!C8+Diffusion!C8+Placeholders!C8"
$wsTotal.Range("D8").Formula = "='Total Load'!D8+Advection!D8+Diffusion!D8+Placeholders!D8"
$wsTotal.Range("E8").Formula = "='Total Load'!E8+Advection!E8+Diffusion!E8+Placeholders!E8"
$wsTotal.Range("F8").Formula = "='Total Load'!F8+Advection!F8+Diffusion!F8+Placeholders!F8"
$wsTotal.Range("G8").Formula = "='Total Load'!G8+Advection!G8+Diffusion!G8+Placeholders!G8"
$wsTotal.Range("H8").Formula = "='Total Load'!H8+Advection!H8+Diffusion!H8+Placeholders!H8"
$wsTotal.Range("I8").Formula = "='Total Load'!I8+Advection!I8+Diffusion!I8+Placeholders!I8"

$wsTotal.Range("C9").Formula = "='Total Load'!C9+Advection!C9+Diffusion!C9+Placeholders!C9"
$wsTotal.Range("D9").Formula = "='Total Load'!D9+Advection!D9+Diffusion!D9+Placeholders!D9"
$wsTotal.Range("E9").Formula = "='Total Load'!E9+Advection!E9+Diffusion!E9+Placeholders!E9"
$wsTotal.Range("F9").Formula = "='Total Load'!F9+Advection!F9+Diffusion!F9+Placeholders!F9"
$wsTotal.Range("G9").Formula = "='Total Load'!G9+Advection!G9+Diffusion!G9+Placeholders!G9"
$wsTotal.Range("H9").Formula = "='Total Load'!H9+Advection!H9+Diffusion!H9+Placeholders!H9"
$wsTotal.Range("I9").Formula = "='Total Load'!I9+Advection!I9+Diffusion!I9+Placeholders!I9"

$wsTotal.Range("R4").Value = 4
$wsTotal.Range("S4").Value = 0
$wsTotal.Range("T4").Value = 1
$wsTotal.Range("U4").Value = 6

$wsTotal.Range("D12").Value = "C_0"
$wsTotal.Range("D13").Value = "C_1"
$wsTotal.Range("D14").Value = "C_2"
$wsTotal.Range("D15").Value = "C_3"
$wsTotal.Range("D16").Value = "C_4"
$wsTotal.Range("D17").Value = "C_5"
$wsTotal.Range("D18").Value = "C_6"

$wsTotal.Range("E12:E18").FormulaArray = "=MMULT(MINVERSE(C3:I9),M3:M9)"

$wsTotal.Rows.Item(12).RowHeight = 45.5
$wsTotal.Rows.Item(13).RowHeight = 45.5
$wsTotal.Rows.Item(14).RowHeight = 45.5
$wsTotal.Rows.Item(15).RowHeight = 45.5
$wsTotal.Rows.Item(16).RowHeight = 45.5
$wsTotal.Rows.Item(17).RowHeight = 45.5
$wsTotal.Rows.Item(18).RowHeight = 45.5

$wsTotal.PageSetup.Orientation = 1

$wsTotal.Activate()
$wsTotal.Range("N8").Select() | Out-Null
